$d = $word.ActiveDocument
$ip = $d.Range(109, 175)
Write-Host "text: [" $ip.Text "]"
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:lang w:val="vi-VN"/></w:rPr><w:tab/><w:t>Hoa don dich vụ của từng công ty: (Sắp xếp giảm dần theo chi phí</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:lang w:val="vi-VN"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:lang w:val="vi-VN"/></w:rPr><w:t>và theo tháng</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:lang w:val="vi-VN"/></w:rPr><w:t>)</w:t></w:r></w:p>
'@
$res = $ip.InsertXML($xml)
Write-Host "Result:" $res
